# Add two new mooring/data-status rows (31 and 32) describing how a large
# source sound file ("LM", mooring M2, BS13_AU_02a) was split into pieces
# that had to be analyzed separately.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 31 and 32 by copying the formatting of an existing, fully
# populated data row (row 2) so the new rows pick up the same number
# formats / styles as the rest of the table (D/E -> 0.00 style, F/G -> date
# style, O -> text style), without allocating any new style entries.
$ws.Range("A2:O2").Copy()
$ws.Range("A31:O31").PasteSpecial(-4122)
$ws.Range("J31:K31").ClearContents()

$ws.Range("A2:O2").Copy()
$ws.Range("A32:O32").PasteSpecial(-4122)
$ws.Range("J32:K32").ClearContents()

# Row 31: LM / M2, 633 files, high-graded via hand mooring selection.
$ws.Range("A31").Value = "LM"
$ws.Range("B31").Value = "M2"
$ws.Range("C31").Value = 633
$ws.Range("D31").Formula = "=L31/C31*100"
$ws.Range("E31").Formula = "=51000/3600"
$ws.Range("F31").Value = 41408
$ws.Range("G31").Value = 41409
$ws.Range("H31").Value = "n"
$ws.Range("I31").Value = "hand (mooring by hand, randomly selected extent) "
$ws.Range("L31").Value = 84
$ws.Range("M31").Value = "n"
$ws.Range("N31").Value = "BS13_AU_02a"
$ws.Range("O31").Value = "38:122"

# Row 32: second piece of the same split source file.
$ws.Range("A32").Value = "LM"
$ws.Range("B32").Value = "M2"
$ws.Range("C32").Value = 633
$ws.Range("D32").Formula = "=L32/C32*100"
$ws.Range("E32").Formula = "=71400/3600"
$ws.Range("F32").Value = 41456
$ws.Range("G32").Value = 41472
$ws.Range("H32").Value = "n"
$ws.Range("I32").Value = "hand (mooring by hand, randomly selected extent) "
$ws.Range("L32").Value = 118
$ws.Range("M32").Value = "n"
$ws.Range("N32").Value = "BS13_AU_02a"
$ws.Range("O32").Value = "510:628"

# Row 32's "Data hours" column (E) keeps the default/general format rather
# than the 0.00 style used elsewhere in column E.
$ws.Range("E32").Style = "Normal"

# Selection / view housekeeping to match the state after the edit: no more
# frozen/scrolled top-left cell, cursor parked just past the new data.
$ws.Range("O33").Select() | Out-Null
